$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -0.2785982940463647
$ws.Range("E2").Value = -0.1790767623711315
$ws.Range("F2").Value = 4.238518120779972
$ws.Range("G2").Value = 33.00453415414271
$ws.Range("H2").Value = 3.9174515734816695
$ws.Range("I2").Value = 61.2
$ws.Range("J2").Value = 0.016999999999999318

$ws.Range("D3").Value = -0.09949436432604675
$ws.Range("E3").Value = -0.2188807929464102
$ws.Range("F3").Value = 1.9501183260035393
$ws.Range("G3").Value = 8.988651242861154
$ws.Range("H3").Value = 2.3003304459706433
$ws.Range("I3").Value = 44.48
$ws.Range("J3").Value = 0.02220000000000027

$ws.Range("D4").Value = 0.0000027836873265529407
$ws.Range("E4").Value = 0.2387885907452136
$ws.Range("F4").Value = 1.631732285121509
$ws.Range("G4").Value = 4.949710440643127
$ws.Range("H4").Value = 1.5276900631062362
$ws.Range("I4").Value = 40.86
$ws.Range("J4").Value = 0.03720000000000482

$ws.Range("D5").Value = 0.019902332514330923
$ws.Range("E5").Value = 0.019905161156993107
$ws.Range("F5").Value = 1.3929421899851324
$ws.Range("G5").Value = 3.880573002656972
$ws.Range("H5").Value = 1.407083050946755
$ws.Range("I5").Value = 38.66
$ws.Range("J5").Value = 0.0596000000000015

$ws.Range("D6").Value = -0.01990019839185063
$ws.Range("E6").Value = 0.09949399861305414
$ws.Range("F6").Value = 0.9153622657766052
$ws.Range("G6").Value = 1.8214944456707816
$ws.Range("H6").Value = 1.0018382940431694
$ws.Range("I6").Value = 36.5
$ws.Range("J6").Value = 0.13060000000000171
